$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoM")

# Update the switch footprint/datasheet/supplier-link for the six Omron
# B3FS tactile-switch rows (Row 28-33, sheet rows 36-41) to the new
# Wuerth 12mm push switch part.
$ws.Range("F36:F41").Value = "SW_PUSH-12mm_Wuerth"
$ws.Range("I36:I41").Value = "https://www.we-online.com/components/products/datasheet/430172043816.pdf"
$ws.Range("J36:J41").Value = "https://www.digikey.ch/de/products/detail/w%C3%BCrth-elektronik/430162043826/9950802"

# Update the component-count summary to reflect the new SMD/THT split.
$ws.Range("F3").Value = "67 (34 SMD/ 31 THT)"
$ws.Range("F4").Value = "67 (34 SMD/ 31 THT)"
